$wb = $excel.ActiveWorkbook

# "总计" is currently the last sheet; insert the new "2022-Q1" sheet right
# before it (Worksheets.Add(Before) mirrors Excel's InsertBefore semantics).
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet = $wb.Worksheets.Item(4)    ; # "2021-Q3" - a donor sheet with the same
                                        ; # header/column-A styling we need to replicate
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# NOTE: $totalSheet was bound by index; inserting a sheet in front of it shifts
# what that index now refers to (it would now resolve to the new sheet).
# Re-fetch "总计" by name so later writes land on the right sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# ---- Fill "2022-Q1" sheet (fund-level detail) ----
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

$rows = @(
    @("470009", "汇添富民营活力混合A", "25.42", "88.81", "4.31", "1.0956", 7),
    @("007355", "汇添富科技创新灵活配置混合A", "16.80", "85.87", "5.18", "0.8702", 3),
    @("009715", "汇添富策略增长两年封闭运作灵活配置混合", "11.81", "88.41", "3.37", "0.3980", 10),
    @("011184", "东方阿尔法招阳混合A", "6.40", "92.60", "5.62", "0.3597", 8),
    @("007356", "汇添富科技创新灵活配置混合C", "2.52", "85.87", "5.18", "0.1305", 3),
    @("013067", "富安达中小盘六个月持有期混合", "2.45", "74.39", "4.58", "0.1122", 1),
    @("011185", "东方阿尔法招阳混合C", "0.08", "92.60", "5.62", "0.0045", 8),
    @("004456", "兴银消费新趋势灵活配置混合", "0.06", "82.19", "5.04", "0.0030", 6),
    @("960014", "汇添富民营活力混合型证券投资基金 O", "0.00", "88.81", "4.31", 0, 7)
)

# Columns B, D:F (and G for all but the last row) hold numeric-looking text such
# as fund codes ("007355") and formatted decimals ("16.80") that must keep their
# original text representation (leading zeros / trailing zeros), so mark them as
# Text before writing - otherwise Excel auto-coerces them to numbers. Column C
# (fund names) always contains non-numeric Chinese text, so it needs no special
# handling.
$newSheet.Range("B2:B10").NumberFormat = "@"
$newSheet.Range("D2:G9").NumberFormat = "@"
$newSheet.Range("D10:F10").NumberFormat = "@"

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $newSheet.Cells.Item($rowNum, 1).Value = $r
    $row = $rows[$r]
    for ($c = 0; $c -lt 6; $c++) {
        $newSheet.Cells.Item($rowNum, 2 + $c).Value = $row[$c]
    }
    $newSheet.Cells.Item($rowNum, 8).Value = $row[6]
}
# G10 is a genuine number (0), unlike G2:G9 which are text - restore default
# number formatting there before writing it.
$newSheet.Range("G10").NumberFormat = "General"
$newSheet.Range("G10").Value = 0

# ---- Match the header / index-column styling used by the other quarter sheets ----
$srcSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$srcSheet.Range("A2:A10").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Insert the new 2022-Q1 summary row at the top of "总计" ----
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 2.97

# Match the bordered index-column style ("s=2") used for the rest of column A.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Renumber the 0-based index column so it stays contiguous after the insert.
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
